$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 734, shifting existing rows 734:775 down to 735:776.
$ws.Rows.Item(734).Insert()

# Populate the newly inserted row with the new data point.
# Column A holds a date-like string; force it to be stored as text (not
# auto-converted to a date serial number) and then clear the formatting
# that was applied so the cell keeps the default style.
$ws.Cells.Item(734, 1).NumberFormat = "@"
$ws.Cells.Item(734, 1).Value = "2026/01/27"
$ws.Cells.Item(734, 1).ClearFormats()

$ws.Cells.Item(734, 2).Value = "火"
$ws.Cells.Item(734, 3).Value = 15
$ws.Cells.Item(734, 4).Value = 201
